# Fruta / hortaliza, semanal
# Insert 4 new daily-price rows for "Terminal Hortofrutícola Agro Chillán - Limón"
# above the existing row 834, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at position 834 (rows 834-888 shift to 838-892).
$ws.Range("A834:A837").EntireRow.Insert()

# Common (constant) column values used throughout this data block.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102003
$categoria = "Limón"
$variedad  = "Sin especificar"

function Set-Row($r, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $origen, $precioKg, $kgUnidad) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 834 44931 "1a amarillo" 120 13000 14000 13500 "$/malla 16 kilos" "Región de O'Higgins" 844 16
Set-Row 835 44931 "1a amarillo" 120 15000 16000 15500 "$/malla 18 kilos" "Provincia de Limarí" 861 18
Set-Row 836 44931 "2a amarillo" 80  12000 12000 12000 "$/malla 16 kilos" "Región de O'Higgins" 750 16
Set-Row 837 44931 "2a amarillo" 80  14000 14000 14000 "$/malla 18 kilos" "Provincia de Limarí" 778 18
